# Fruta / hortaliza, semanal
#
# Insert a new week of price data (row group for 2022-09-22, serial 44826)
# for "Packham's Triumph" (Especial/Primera/Segunda) right above the
# existing 2021-08-05 (serial 44413) row group. This pushes every
# subsequent data row down by 3 (old rows 790-843 become 793-846).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 3 blank rows at 790:792 - this shifts the old 790:792 (and
#    everything below) down to 793:795 .. 846.
$ws.Rows("790:792").Insert()

# 2) Seed the newly-inserted blank rows with a copy of what is now at
#    793:795 (the original 790:792 content), so every column (product,
#    variety, quality, units, region, etc.) starts out identical.
$ws.Range("A790:T792").Value = $ws.Range("A793:T795").Value()

# 3) Overwrite the new week's date + volume/price figures.
#    Row 790: Packham's Triumph / Especial
$ws.Range("D790").Value = 44826
$ws.Range("M790").Value = 20
$ws.Range("N790").Value = 310000
$ws.Range("O790").Value = 320000
$ws.Range("P790").Value = 315000
$ws.Range("S790").Value = 700

#    Row 791: Packham's Triumph / Primera
$ws.Range("D791").Value = 44826
$ws.Range("M791").Value = 16
$ws.Range("N791").Value = 290000
$ws.Range("O791").Value = 300000
$ws.Range("P791").Value = 295000
$ws.Range("S791").Value = 656

#    Row 792: Packham's Triumph / Segunda
$ws.Range("D792").Value = 44826
$ws.Range("M792").Value = 14
$ws.Range("N792").Value = 255000
$ws.Range("O792").Value = 260000
$ws.Range("P792").Value = 257500
$ws.Range("S792").Value = 572
